$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 9025024
$ws.Range("I80").Value = 641.2857
$ws.Range("K80").Value = 1923.8571
$ws.Range("M80").Value = -925.8571000000002

$ws.Range("H83").Value = 9025024
$ws.Range("I83").Value = 641.2857
$ws.Range("K83").Value = 5771.571300000001
$ws.Range("M83").Value = -779.5713000000005

$ws.Range("H92").Value = 111111944
$ws.Range("I92").Value = 111111944
$ws.Range("K92").Value = 111111944
$ws.Range("M92").Value = -111110696

$ws.Range("H109").Value = 30634.5
$ws.Range("J109").Value = 30634.5
$ws.Range("L109").Value = 30634.5
$ws.Range("N109").Value = -33408.5

$ws.Range("H113").Value = 38465708
$ws.Range("I113").Value = 83336860
$ws.Range("J113").Value = 4714.0713
$ws.Range("K113").Value = 83336860
$ws.Range("L113").Value = 4714.0713
$ws.Range("M113").Value = -83333606
$ws.Range("N113").Value = -11222.0713

$ws.Range("H137").Value = 31667.908
$ws.Range("I137").Value = 1298.24
$ws.Range("J137").Value = 126573.125
$ws.Range("K137").Value = 3894.72
$ws.Range("L137").Value = 379719.375
$ws.Range("M137").Value = -1344.72
$ws.Range("N137").Value = -384819.375

$ws.Range("H138").Value = 1767.0116
$ws.Range("I138").Value = 743.14813
$ws.Range("J138").Value = 2235.5593
$ws.Range("K138").Value = 2229.44439
$ws.Range("L138").Value = 6706.6779
$ws.Range("M138").Value = 2910.55561
$ws.Range("N138").Value = -16986.6779

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1378.9474
$ws.Range("I107").Value = 1185.5714
$ws.Range("J107").Value = 1920.4
$ws.Range("K107").Value = 1185.5714
$ws.Range("L107").Value = 1920.4
$ws.Range("M107").Value = 734.4286
$ws.Range("N107").Value = -5760.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 422.42856
$ws.Range("I5").Value = 101.5
$ws.Range("J5").Value = 550.8
$ws.Range("K5").Value = 101.5
$ws.Range("L5").Value = 550.8
$ws.Range("M5").Value = 10.5
$ws.Range("N5").Value = -774.8

$ws.Range("H31").Value = 10205.475
$ws.Range("I31").Value = 20077.295
$ws.Range("J31").Value = 2908.913
$ws.Range("K31").Value = 20077.295
$ws.Range("L31").Value = 2908.913
$ws.Range("M31").Value = -19782.295
$ws.Range("N31").Value = -3498.913

$ws.Range("H34").Value = 10205.475
$ws.Range("I34").Value = 20077.295
$ws.Range("J34").Value = 2908.913
$ws.Range("K34").Value = 20077.295
$ws.Range("L34").Value = 2908.913
$ws.Range("M34").Value = -19875.295
$ws.Range("N34").Value = -3312.913

$ws.Range("H122").Value = 1086.561
$ws.Range("I122").Value = 1093.1765
$ws.Range("K122").Value = 3279.5295
$ws.Range("M122").Value = -829.5295000000001

$ws.Range("H134").Value = 953.13043
$ws.Range("I134").Value = 896.1
$ws.Range("K134").Value = 2688.3
$ws.Range("M134").Value = -153.3000000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3406.5
$ws.Range("I75").Value = 1813
$ws.Range("J75").Value = 5000
$ws.Range("K75").Value = 5439
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = -4441
$ws.Range("N75").Value = -16996

$ws.Range("H78").Value = 3406.5
$ws.Range("I78").Value = 1813
$ws.Range("J78").Value = 5000
$ws.Range("K78").Value = 16317
$ws.Range("L78").Value = 45000
$ws.Range("M78").Value = -11325
$ws.Range("N78").Value = -54984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2978887.8
$ws.Range("I11").Value = 3215713.8
$ws.Range("J11").Value = 2149997
$ws.Range("K11").Value = 3215713.8
$ws.Range("L11").Value = 2149997
$ws.Range("M11").Value = -3215574.8
$ws.Range("N11").Value = -2150275

$ws.Range("H12").Value = 6556250
$ws.Range("J12").Value = 6000000
$ws.Range("L12").Value = 6000000
$ws.Range("N12").Value = -6000280

$ws.Range("H14").Value = 3673728.8
$ws.Range("I14").Value = 4041051.5
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 4041051.5
$ws.Range("L14").Value = 500
$ws.Range("M14").Value = -4040883.5
$ws.Range("N14").Value = -836

$ws.Range("H80").Value = 3149.9092
$ws.Range("I80").Value = 2730.5
$ws.Range("J80").Value = 3499.4167
$ws.Range("K80").Value = 2730.5
$ws.Range("L80").Value = 3499.4167
$ws.Range("M80").Value = -1732.5
$ws.Range("N80").Value = -5495.4167

$ws.Range("H83").Value = 3149.9092
$ws.Range("I83").Value = 2730.5
$ws.Range("J83").Value = 3499.4167
$ws.Range("K83").Value = 13652.5
$ws.Range("L83").Value = 17497.0835
$ws.Range("M83").Value = -8660.5
$ws.Range("N83").Value = -27481.0835

$ws.Range("H95").Value = 4833.3335
$ws.Range("J95").Value = 4833.3335
$ws.Range("L95").Value = 4833.3335
$ws.Range("N95").Value = -10325.3335

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4952.381
$ws.Range("I7").Value = 2910
$ws.Range("K7").Value = 2910
$ws.Range("M7").Value = -2798

$ws.Range("H14").Value = 2940
$ws.Range("J14").Value = 2940
$ws.Range("L14").Value = 2940
$ws.Range("N14").Value = -3284

$ws.Range("H25").Value = 4500
$ws.Range("J25").Value = 4500
$ws.Range("L25").Value = 4500
$ws.Range("N25").Value = -4960

$ws.Range("H40").Value = 6713.25
$ws.Range("I40").Value = 5166.5
$ws.Range("J40").Value = 8260
$ws.Range("K40").Value = 5166.5
$ws.Range("L40").Value = 8260
$ws.Range("M40").Value = -5030.5
$ws.Range("N40").Value = -8532

$ws.Range("H61").Value = 4975.9165
$ws.Range("J61").Value = 6971.273
$ws.Range("L61").Value = 6971.273
$ws.Range("N61").Value = -7375.273

$ws.Range("H74").Value = 31722
$ws.Range("J74").Value = 48000
$ws.Range("L74").Value = 48000
$ws.Range("N74").Value = -49996

$ws.Range("H77").Value = 31722
$ws.Range("J77").Value = 48000
$ws.Range("L77").Value = 144000
$ws.Range("N77").Value = -153984

$ws.Range("H82").Value = 2155.8
$ws.Range("J82").Value = 1992.7142
$ws.Range("L82").Value = 1992.7142
$ws.Range("N82").Value = -2714.7142

$ws.Range("H85").Value = 2155.8
$ws.Range("J85").Value = 1992.7142
$ws.Range("L85").Value = 1992.7142
$ws.Range("N85").Value = -4488.7142

$ws.Range("H93").Value = 1475.5
$ws.Range("I93").Value = 1475.5
$ws.Range("K93").Value = 1475.5
$ws.Range("M93").Value = -227.5

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H113").Value = 4975.9165
$ws.Range("J113").Value = 6971.273
$ws.Range("L113").Value = 6971.273
$ws.Range("N113").Value = -11311.273

$ws.Range("H126").Value = 4952.381
$ws.Range("I126").Value = 2910
$ws.Range("K126").Value = 8730
$ws.Range("M126").Value = -6260

$ws.Range("H132").Value = 1875.48
$ws.Range("I132").Value = 1370.4706
$ws.Range("J132").Value = 2948.625
$ws.Range("K132").Value = 4111.4118
$ws.Range("L132").Value = 8845.875
$ws.Range("M132").Value = -1581.4118
$ws.Range("N132").Value = -13905.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 5000100
$ws.Range("I17").Value = 5000100
$ws.Range("K17").Value = 5000100
$ws.Range("M17").Value = -4999928

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H80").Value = 34650
$ws.Range("J80").Value = 34650
$ws.Range("L80").Value = 34650
$ws.Range("N80").Value = -36646

$ws.Range("H81").Value = 100001720
$ws.Range("I81").Value = 2050.125
$ws.Range("J81").Value = 500000420
$ws.Range("K81").Value = 4100.25
$ws.Range("L81").Value = 1000000840
$ws.Range("M81").Value = -3039.25
$ws.Range("N81").Value = -1000002962

$ws.Range("H83").Value = 34650
$ws.Range("J83").Value = 34650
$ws.Range("L83").Value = 103950
$ws.Range("N83").Value = -113934

$ws.Range("H84").Value = 100001720
$ws.Range("I84").Value = 2050.125
$ws.Range("J84").Value = 500000420
$ws.Range("K84").Value = 20501.25
$ws.Range("L84").Value = 5000004200
$ws.Range("M84").Value = -15197.25
$ws.Range("N84").Value = -5000014808

$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344

$ws.Range("H113").Value = 1352854.6
$ws.Range("I113").Value = 2188.3076
$ws.Range("K113").Value = 6564.9228
$ws.Range("M113").Value = -4394.9228
